$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "NoRedistribution" column (P) ---------------------------
$ws.Range("P1").Value = "NoRedistribution"

# Rows 2-16 correspond to A2:A16 (licenses + the special "Owner" row at 15).
# Only the "Owner" row is allowed to redistribute privately-owned code, so
# every license row is "Excluded" and the Owner row is "Included".
for ($r = 2; $r -le 16; $r++) {
    if ($r -eq 15) {
        $ws.Cells.Item($r, 16).Value = "Included"
    } else {
        $ws.Cells.Item($r, 16).Value = "Excluded"
    }
}

# --- Conditional formatting for the new column ----------------------------
# Same "red" highlight style already used by the other Included/Excluded
# columns, applied to the new P2:P16 range.
$fcP = $ws.Range("P2:P16").FormatConditions.Add(1, 3, '"Included"')
$fcP.Font.Color = 393372
$fcP.Interior.Color = 13551615

# The newly-added rule takes the highest precedence (priority 1), same as
# Excel does when a conditional format is added through the UI; bump the
# pre-existing rules down to make room for it.
$fcL = $ws.Range("L2:O16").FormatConditions.Item(1)
$fcG = $ws.Range("G2:K16").FormatConditions.Item(1)
$fcB = $ws.Range("B2:F16").FormatConditions.Item(1)

$fcP.SetFirstPriority()
$fcL.Priority = 2
$fcG.Priority = 3
$fcB.Priority = 4

# --- Match the resulting selection left behind by the edit ---------------
$ws.Range("P16").Select()
